$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 749.4
$ws.Range("J53").Value = 600
$ws.Range("L53").Value = 600
$ws.Range("N53").Value = -1874
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H62").Value = 3502
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3502
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 1911.5
$ws.Range("I94").Value = 1713.1428
$ws.Range("K94").Value = 1713.1428
$ws.Range("M94").Value = -1262.1428
$ws.Range("H98").Value = 5406.522
$ws.Range("I98").Value = 5073
$ws.Range("J98").Value = 7630
$ws.Range("K98").Value = 5073
$ws.Range("L98").Value = 7630
$ws.Range("M98").Value = -3575
$ws.Range("N98").Value = -10626
$ws.Range("H112").Value = 406642.6
$ws.Range("I112").Value = 3192.5
$ws.Range("J112").Value = 476807.8
$ws.Range("K112").Value = 9577.5
$ws.Range("L112").Value = 1430423.4
$ws.Range("M112").Value = -8469.5
$ws.Range("N112").Value = -1432639.4
$ws.Range("H122").Value = 5406.522
$ws.Range("I122").Value = 5073
$ws.Range("J122").Value = 7630
$ws.Range("K122").Value = 15219
$ws.Range("L122").Value = 22890
$ws.Range("M122").Value = -12769
$ws.Range("N122").Value = -27790
$ws.Range("H132").Value = 2598.7144
$ws.Range("I132").Value = 2260.6924
$ws.Range("J132").Value = 6993
$ws.Range("K132").Value = 6782.0772
$ws.Range("L132").Value = 20979
$ws.Range("M132").Value = -4252.0772
$ws.Range("N132").Value = -26039
$ws.Range("H137").Value = 2441.125
$ws.Range("I137").Value = 1517.8
$ws.Range("K137").Value = 4553.4
$ws.Range("M137").Value = -2003.4
$ws.Range("H138").Value = 2947.5
$ws.Range("I138").Value = 1615.5
$ws.Range("J138").Value = 3063.3262
$ws.Range("K138").Value = 4846.5
$ws.Range("L138").Value = 9189.9786
$ws.Range("M138").Value = 293.5
$ws.Range("N138").Value = -19469.9786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 51357.71
$ws.Range("I32").Value = 11951.184
$ws.Range("K32").Value = 11951.184
$ws.Range("M32").Value = -11664.184
$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -19685
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3869.5366
$ws.Range("I134").Value = 4028.342
$ws.Range("K134").Value = 12085.026
$ws.Range("M134").Value = -9550.026

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6063018.5
$ws.Range("J31").Value = 2541.6667
$ws.Range("L31").Value = 2541.6667
$ws.Range("N31").Value = -3131.6667
$ws.Range("H34").Value = 6063018.5
$ws.Range("J34").Value = 2541.6667
$ws.Range("L34").Value = 2541.6667
$ws.Range("N34").Value = -2945.6667
$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 8000
$ws.Range("K51").Value = 8000
$ws.Range("M51").Value = -7264
$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 8000
$ws.Range("K61").Value = 8000
$ws.Range("M61").Value = -7652
$ws.Range("H62").Value = 3599.3845
$ws.Range("I62").Value = 3818.182
$ws.Range("J62").Value = 2396
$ws.Range("K62").Value = 3818.182
$ws.Range("L62").Value = 2396
$ws.Range("M62").Value = -3194.182
$ws.Range("N62").Value = -3644
$ws.Range("H65").Value = 3599.3845
$ws.Range("I65").Value = 3818.182
$ws.Range("J65").Value = 2396
$ws.Range("K65").Value = 19090.91
$ws.Range("L65").Value = 11980
$ws.Range("M65").Value = -15970.91
$ws.Range("N65").Value = -18220
$ws.Range("H99").Value = 3183.476
$ws.Range("I99").Value = 3259.8572
$ws.Range("K99").Value = 3259.8572
$ws.Range("M99").Value = -1761.8572
$ws.Range("H126").Value = 3183.476
$ws.Range("I126").Value = 3259.8572
$ws.Range("K126").Value = 9779.571599999999
$ws.Range("M126").Value = -7309.571599999999
$ws.Range("H132").Value = 4330.6
$ws.Range("I132").Value = 4126.2
$ws.Range("K132").Value = 12378.6
$ws.Range("M132").Value = -9848.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3145.1924
$ws.Range("J34").Value = 3432.3044
$ws.Range("L34").Value = 10296.9132
$ws.Range("N34").Value = -10464.9132
$ws.Range("H47").Value = 3208.2144
$ws.Range("I47").Value = 102.5
$ws.Range("J47").Value = 5537.5
$ws.Range("K47").Value = 307.5
$ws.Range("L47").Value = 16612.5
$ws.Range("M47").Value = 123.5
$ws.Range("N47").Value = -17474.5
$ws.Range("H48").Value = 7200
$ws.Range("J48").Value = 11666.667
$ws.Range("L48").Value = 35000.001
$ws.Range("N48").Value = -35500.001
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H50").Value = 743.7
$ws.Range("I50").Value = 716.3333
$ws.Range("J50").Value = 784.75
$ws.Range("K50").Value = 2148.9999
$ws.Range("L50").Value = 2354.25
$ws.Range("M50").Value = -1667.9999
$ws.Range("N50").Value = -3316.25
$ws.Range("H53").Value = 743.7
$ws.Range("I53").Value = 716.3333
$ws.Range("J53").Value = 784.75
$ws.Range("K53").Value = 2148.9999
$ws.Range("L53").Value = 2354.25
$ws.Range("M53").Value = -1667.9999
$ws.Range("N53").Value = -3316.25
$ws.Range("H54").Value = 12600
$ws.Range("I54").Value = 15000
$ws.Range("J54").Value = 12000
$ws.Range("K54").Value = 45000
$ws.Range("L54").Value = 36000
$ws.Range("M54").Value = -44441
$ws.Range("N54").Value = -37118
$ws.Range("H55").Value = 5250
$ws.Range("J55").Value = 5250
$ws.Range("L55").Value = 15750
$ws.Range("N55").Value = -16104
$ws.Range("H56").Value = 5739.5713
$ws.Range("I56").Value = 5739.5713
$ws.Range("K56").Value = 5739.5713
$ws.Range("M56").Value = -5209.5713
$ws.Range("H59").Value = 450
$ws.Range("I59").Value = 450
$ws.Range("K59").Value = 1350
$ws.Range("M59").Value = -810
$ws.Range("H60").Value = 2051.875
$ws.Range("J60").Value = 156
$ws.Range("L60").Value = 468
$ws.Range("N60").Value = -970
$ws.Range("H61").Value = 99
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H97").Value = 10307.5
$ws.Range("I97").Value = 357.33334
$ws.Range("K97").Value = 1072.00002
$ws.Range("M97").Value = -576.0000199999999
$ws.Range("H113").Value = 1150.9231
$ws.Range("J113").Value = 1150.9231
$ws.Range("L113").Value = 3452.7693
$ws.Range("N113").Value = -7792.7693
$ws.Range("H114").Value = 1501.909
$ws.Range("I114").Value = 1115
$ws.Range("J114").Value = 1587.8889
$ws.Range("K114").Value = 3345
$ws.Range("L114").Value = 4763.6667
$ws.Range("M114").Value = -91
$ws.Range("N114").Value = -11271.6667
$ws.Range("H131").Value = 1524.7307
$ws.Range("I131").Value = 1396.4166
$ws.Range("J131").Value = 1634.7142
$ws.Range("K131").Value = 4189.2498
$ws.Range("L131").Value = 4904.142599999999
$ws.Range("M131").Value = 850.7502000000004
$ws.Range("N131").Value = -14984.1426
$ws.Range("H132").Value = 7156.143
$ws.Range("I132").Value = 2783.1667
$ws.Range("J132").Value = 10435.875
$ws.Range("K132").Value = 25048.5003
$ws.Range("L132").Value = 93922.875
$ws.Range("M132").Value = -22518.5003
$ws.Range("N132").Value = -98982.875
$ws.Range("H139").Value = 2276.0833
$ws.Range("I139").Value = 2501.75
$ws.Range("K139").Value = 7505.25
$ws.Range("M139").Value = -2365.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 32497
$ws.Range("J117").Value = 32497
$ws.Range("L117").Value = 32497
$ws.Range("N117").Value = -39381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 24502.8
$ws.Range("H20").Value = 888400
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H93").Value = 2752.9714
$ws.Range("I93").Value = 1987
$ws.Range("K93").Value = 1987
$ws.Range("M93").Value = -739
$ws.Range("H136").Value = 38466972
$ws.Range("I136").Value = 5217.913
$ws.Range("K136").Value = 15653.739
$ws.Range("M136").Value = -13103.739

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 36434.75
$ws.Range("J101").Value = 36434.75
$ws.Range("L101").Value = 36434.75
$ws.Range("N101").Value = -42924.75
